# Actualización automática 2025-10-20 14:30:09
# Registers a new sale of 1758.38 for ALMEIDA CUATIN JHONATHANN CARLOS /
# ZAMBRANO CEDEÑO MARJORIE XIOMARA (PORCELANATO) in october, and propagates
# the totals across the three report sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M35").Value = 1758.38
$wsGrupo.Range("M36").Value = "8 de 34"

# --- Sheet "VENTA MENSUAL" -------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F35").Value = 1758.38
$wsMensual.Range("F36").Value = 18748.37

# --- Sheet "CUMPLIMIENTO MENSUAL" ------------------------------------------
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D12").Value = 16871.89
$wsCumplimiento.Range("E12").Value = 4829.380000000001
$wsCumplimiento.Range("F12").Value = 0.7774609504420709

$wsCumplimiento.Range("D14").Value = 18748.37
$wsCumplimiento.Range("E14").Value = 17837.19723718183
$wsCumplimiento.Range("F14").Value = 0.5124526258799146
